$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format first so numeric-looking strings
# (e.g. "1.004", "291.71") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.393.27'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.571.96'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '1.004'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '291.71'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = '0.3761'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').Value = '50.07'
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  +0.49%  '
$ws.Range('D11').Value = '1.153'
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '6.014'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').Value = '6.959'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').Value = '1.572.41'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '0.00001131'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '90.01'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '16.78'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('D22').Value = '6.219'
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = '11.98'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = '22.394.28'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '2.404'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').Value = '2.672'
$ws.Range('E26').Value = '  -10.37%  '
$ws.Range('D27').Value = '20.19'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').Value = '147.26'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '5.040'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').Value = '126.49'
$ws.Range('E30').Value = '  +0.87%  '
$ws.Range('D31').Value = '1.754.83'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '6.148'
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('D33').Value = '2.002'
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('D34').Value = '0.9833'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('D35').Value = '9.992'
$ws.Range('E35').Value = '  -3.16%  '
$ws.Range('D36').Value = '0.08487'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '0.02538'
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').Value = '1.375'
$ws.Range('E38').Value = '  +10.45%  '
$ws.Range('D39').Value = '0.2313'
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('D40').Value = '0.06576'
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').Value = '5.412'
$ws.Range('E41').Value = '  -2.52%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '11.46'
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.6389'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '14.04'
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '0.5975'
$ws.Range('E47').Value = '  -0.64%  '
$ws.Range('D48').Value = '1.296'
$ws.Range('D49').Value = '2.091'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('D50').Value = '125.41'
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').Value = '0.07332'
$ws.Range('E51').Value = '  +0.66%  '

# Restore the normal style on column D so no residual number format remains
# on cells (matches original formatting -- values stay as text).
$ws.Range("D2:D51").Style = "Normal"
